$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Result")

$ws.Range("A1").Value = "ISBN"
$ws.Range("B1").Value = "Tittel"
$ws.Range("C1").Value = "Forlag"
